# Refresh crypto price/volume figures (cron-style data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Columns D/E hold plain text (prices with "." thousands separators,
    # percentages, etc.) -- force text interpretation so Excel does not
    # auto-coerce numeric-looking strings into numbers, then restore the
    # default "Normal" style so no stray formatting is left behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '61.030.57'
Set-TextValue $ws.Range("E2") '  +1.47%  '
Set-TextValue $ws.Range("D3") '2.350.97'
Set-TextValue $ws.Range("E3") '  +0.77%  '
Set-TextValue $ws.Range("E4") '  -0.04%  '
Set-TextValue $ws.Range("D5") '557.03'
Set-TextValue $ws.Range("E5") '  +2.42%  '
Set-TextValue $ws.Range("D6") '132.27'
Set-TextValue $ws.Range("E6") '  +0.73%  '
Set-TextValue $ws.Range("D7") '0.999'
Set-TextValue $ws.Range("E7") '  -0.04%  '
Set-TextValue $ws.Range("D8") '0.583'
Set-TextValue $ws.Range("E8") '  -0.28%  '
Set-TextValue $ws.Range("D9") '2.348.37'
Set-TextValue $ws.Range("E9") '  +0.79%  '
Set-TextValue $ws.Range("D10") '0.103'
Set-TextValue $ws.Range("E10") '  +1.69%  '
Set-TextValue $ws.Range("D11") '5.63'
Set-TextValue $ws.Range("E11") '  +2.34%  '
Set-TextValue $ws.Range("E12") '  -0.41%  '
Set-TextValue $ws.Range("D13") '0.340'
Set-TextValue $ws.Range("E13") '  +1.77%  '
Set-TextValue $ws.Range("D14") '24.14'
Set-TextValue $ws.Range("E14") '  +1.79%  '
Set-TextValue $ws.Range("D15") '2.768.34'
Set-TextValue $ws.Range("E15") '  +0.77%  '
Set-TextValue $ws.Range("D16") '60.856.89'
Set-TextValue $ws.Range("E16") '  +1.29%  '
Set-TextValue $ws.Range("E17") '  +2.07%  '
Set-TextValue $ws.Range("D18") '2.333.64'
Set-TextValue $ws.Range("E18") '  +0.60%  '
Set-TextValue $ws.Range("D19") '10.74'
Set-TextValue $ws.Range("E19") '  +1.39%  '
Set-TextValue $ws.Range("D20") '4.13'
Set-TextValue $ws.Range("E20") '  -0.51%  '
Set-TextValue $ws.Range("D21") '316.69'
Set-TextValue $ws.Range("E21") '  +1.00%  '
Set-TextValue $ws.Range("D22") '6.67'
Set-TextValue $ws.Range("E22") '  -1.27%  '
Set-TextValue $ws.Range("E23") '  +0.38%  '
Set-TextValue $ws.Range("D24") '64.38'
Set-TextValue $ws.Range("E24") '  +1.31%  '
Set-TextValue $ws.Range("D25") '0.172'
Set-TextValue $ws.Range("E25") '  +0.38%  '
Set-TextValue $ws.Range("E26") '  -0.16%  '
Set-TextValue $ws.Range("D27") '8.04'
Set-TextValue $ws.Range("E27") '  +1.73%  '
Set-TextValue $ws.Range("D28") '1.43'
Set-TextValue $ws.Range("E28") '  +5.39%  '
Set-TextValue $ws.Range("D29") '1.29'
Set-TextValue $ws.Range("E29") '  +10.19%  '
Set-TextValue $ws.Range("E30") '  +0.69%  '
Set-TextValue $ws.Range("D31") '171.33'
Set-TextValue $ws.Range("E31") '  -0.32%  '
Set-TextValue $ws.Range("D32") '0.0₃0740'
Set-TextValue $ws.Range("E32") '  +1.67%  '
Set-TextValue $ws.Range("D33") '6.15'
Set-TextValue $ws.Range("E33") '  +3.66%  '
Set-TextValue $ws.Range("E34") '  -0.10%  '
Set-TextValue $ws.Range("D35") '0.388'
Set-TextValue $ws.Range("E35") '  +1.75%  '
Set-TextValue $ws.Range("D36") '18.14'
Set-TextValue $ws.Range("E36") '  +0.87%  '
Set-TextValue $ws.Range("E38") '  -0.08%  '
Set-TextValue $ws.Range("D39") '4.18'
Set-TextValue $ws.Range("E39") '  +1.45%  '
Set-TextValue $ws.Range("D40") '336.77'
Set-TextValue $ws.Range("E40") '  +5.29%  '
Set-TextValue $ws.Range("D41") '1.55'
Set-TextValue $ws.Range("E41") '  +1.98%  '
Set-TextValue $ws.Range("D42") '38.16'
Set-TextValue $ws.Range("E42") '  +0.28%  '
Set-TextValue $ws.Range("D43") '139.96'
Set-TextValue $ws.Range("E43") '  -0.23%  '
Set-TextValue $ws.Range("D44") '3.55'
Set-TextValue $ws.Range("E44") '  +2.71%  '
Set-TextValue $ws.Range("D45") '0.0953'
Set-TextValue $ws.Range("E45") '  +0.91%  '
Set-TextValue $ws.Range("D46") '19.49'
Set-TextValue $ws.Range("E46") '  +0.23%  '
Set-TextValue $ws.Range("D47") '0.575'
Set-TextValue $ws.Range("E47") '  +2.74%  '
Set-TextValue $ws.Range("D48") '0.0503'
Set-TextValue $ws.Range("E48") '  +1.24%  '
Set-TextValue $ws.Range("D49") '0.0₆0230'
Set-TextValue $ws.Range("E49") '  +8.77%  '
Set-TextValue $ws.Range("E50") '  +2.59%  '
Set-TextValue $ws.Range("D51") '17.29'
Set-TextValue $ws.Range("E51") '  +2.84%  '
